# "remove column from alcohol data"
# The measurement sheet (Sheet1) has an extra column M that duplicates/
# shadows column N's data. Delete column M so the old column N shifts
# left and becomes the new column M (dimension shrinks from N to M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Delete the entire column M (rows 1-119 hold data; EntireColumn makes sure
# the whole column is removed and everything to the right shifts left).
$ws.Range("M1:M119").EntireColumn.Delete()

# Leave the selection on the (now last) column, matching the saved view.
$ws.Range("M1").Select()
